$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repull/recalculation of data
$ws.Cells.Item(2, 6).Value = -4
$ws.Cells.Item(5, 6).Value = 2
$ws.Cells.Item(6, 6).Value = -2
$ws.Cells.Item(10, 6).Value = -1
$ws.Cells.Item(11, 6).Value = 3
$ws.Cells.Item(14, 6).Value = -3
$ws.Cells.Item(15, 6).Value = -4
$ws.Cells.Item(16, 6).Value = -11
$ws.Cells.Item(18, 6).Value = 3
$ws.Cells.Item(19, 6).Value = -5
$ws.Cells.Item(21, 6).Value = -2
$ws.Cells.Item(23, 6).Value = -3
$ws.Cells.Item(30, 6).Value = -7
$ws.Cells.Item(35, 6).Value = -2
$ws.Cells.Item(38, 6).Value = 3
$ws.Cells.Item(41, 6).Value = -3
$ws.Cells.Item(44, 6).Value = 4
$ws.Cells.Item(47, 6).Value = -1
$ws.Cells.Item(49, 6).Value = 1
$ws.Cells.Item(50, 6).Value = -6
$ws.Cells.Item(51, 6).Value = 2
$ws.Cells.Item(52, 6).Value = -3
$ws.Cells.Item(54, 6).Value = -4
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(57, 6).Value = 2
$ws.Cells.Item(58, 6).Value = -3
